# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
# Replaces the previous "estado de cuenta" (account statement) detail rows
# (B16:G36) with the refreshed data for the three workers, grouped by
# worker and ordered by period (most recent period first), with an
# updated "Valor Mora" total per worker.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$workers = @(
    @{ Doc = "9297354";    Name = "EDINSON JOSE RAMOS ZAMBRANO";    Mora = 1000000 },
    @{ Doc = "1050953591"; Name = "YAIRA MARGARITA MUÑIZ GAMARRA";  Mora = 1300000 },
    @{ Doc = "1143385667"; Name = "DANILO HERNANDEZ RIVERA";        Mora = 0 }
)

# Periods, most recent first - each worker has one row per period.
$periods = @("2410", "2409", "2408", "2407", "2406", "2405", "2404")

$row = 16
foreach ($worker in $workers) {
    $isFirstPeriodForWorker = $true
    foreach ($period in $periods) {
        $ws.Range("B$row").Value = "CC"
        $ws.Range("C$row").Value = $worker.Doc
        $ws.Range("D$row").Value = $worker.Name
        $ws.Range("E$row").Value = $period

        if ($isFirstPeriodForWorker) {
            $ws.Range("F$row").Value = 15600
            $isFirstPeriodForWorker = $false
        } else {
            $ws.Range("F$row").Value = 52000
        }

        $ws.Range("G$row").Value = $worker.Mora

        $row = $row + 1
    }
}
